# Auto-generated Excel COM-interop edit script
# Applies numeric corrections to the Leve profit tables across all 8 job sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the scheduled-runner price refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4778.6665
$ws.Range("I64").Value = 4166.5
$ws.Range("J64").Value = 6003
$ws.Range("K64").Value = 4166.5
$ws.Range("L64").Value = 6003
$ws.Range("M64").Value = -3918.5
$ws.Range("N64").Value = -6499
$ws.Range("H67").Value = 4778.6665
$ws.Range("I67").Value = 4166.5
$ws.Range("J67").Value = 6003
$ws.Range("K67").Value = 4166.5
$ws.Range("L67").Value = 6003
$ws.Range("M67").Value = -3308.5
$ws.Range("N67").Value = -7719
$ws.Range("H76").Value = 3091
$ws.Range("I76").Value = 2917.25
$ws.Range("J76").Value = 3322.6667
$ws.Range("K76").Value = 2917.25
$ws.Range("L76").Value = 3322.6667
$ws.Range("M76").Value = -2602.25
$ws.Range("N76").Value = -3952.6667
$ws.Range("H79").Value = 3091
$ws.Range("I79").Value = 2917.25
$ws.Range("J79").Value = 3322.6667
$ws.Range("K79").Value = 2917.25
$ws.Range("L79").Value = 3322.6667
$ws.Range("M79").Value = -1825.25
$ws.Range("N79").Value = -5506.6667
$ws.Range("H132").Value = 1825
$ws.Range("I132").Value = 1842.8572
$ws.Range("K132").Value = 5528.571599999999
$ws.Range("M132").Value = -2998.571599999999
$ws.Range("H141").Value = 2681
$ws.Range("J141").Value = 9998.5
$ws.Range("L141").Value = 29995.5
$ws.Range("N141").Value = -40355.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2100
$ws.Range("J45").Value = 2500
$ws.Range("L45").Value = 2500
$ws.Range("N45").Value = -3254
$ws.Range("H63").Value = 1311.5
$ws.Range("I63").Value = 1138.9
$ws.Range("K63").Value = 1138.9
$ws.Range("M63").Value = -452.9000000000001
$ws.Range("H66").Value = 1311.5
$ws.Range("I66").Value = 1138.9
$ws.Range("K66").Value = 5694.5
$ws.Range("M66").Value = -2262.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 491.66666
$ws.Range("I20").Value = 513.2857
$ws.Range("J20").Value = 416
$ws.Range("K20").Value = 513.2857
$ws.Range("L20").Value = 416
$ws.Range("M20").Value = -266.2857
$ws.Range("N20").Value = -910
$ws.Range("H94").Value = 1949.0625
$ws.Range("I94").Value = 2057.0833
$ws.Range("K94").Value = 2057.0833
$ws.Range("M94").Value = -1606.0833
$ws.Range("H105").Value = 1425
$ws.Range("I105").Value = 350
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 350
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = 1397
$ws.Range("N105").Value = -5994

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2399.6667
$ws.Range("I58").Value = 2399.5
$ws.Range("J58").Value = 2400
$ws.Range("K58").Value = 2399.5
$ws.Range("L58").Value = 2400
$ws.Range("M58").Value = -2196.5
$ws.Range("N58").Value = -2806
$ws.Range("H76").Value = 5000
$ws.Range("I76").Value = 5000
$ws.Range("K76").Value = 5000
$ws.Range("M76").Value = -4685
$ws.Range("H79").Value = 5000
$ws.Range("I79").Value = 5000
$ws.Range("K79").Value = 5000
$ws.Range("M79").Value = -3908
$ws.Range("H107").Value = 1401
$ws.Range("I107").Value = 1614
$ws.Range("J107").Value = 975
$ws.Range("K107").Value = 1614
$ws.Range("L107").Value = 975
$ws.Range("M107").Value = 306
$ws.Range("N107").Value = -4815
$ws.Range("H132").Value = 1120
$ws.Range("I132").Value = 1120
$ws.Range("K132").Value = 3360
$ws.Range("M132").Value = -830
$ws.Range("H136").Value = 2399.6667
$ws.Range("I136").Value = 2399.5
$ws.Range("J136").Value = 2400
$ws.Range("K136").Value = 7198.5
$ws.Range("L136").Value = 7200
$ws.Range("M136").Value = -4648.5
$ws.Range("N136").Value = -12300

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 200
$ws.Range("I2").Value = 15.5
$ws.Range("K2").Value = 93
$ws.Range("M2").Value = 20
$ws.Range("H38").Value = 359.8
$ws.Range("I38").Value = 74.75
$ws.Range("K38").Value = 224.25
$ws.Range("M38").Value = 122.75
$ws.Range("H117").Value = 974.75
$ws.Range("I117").Value = 849.5
$ws.Range("K117").Value = 2548.5
$ws.Range("M117").Value = 893.5
$ws.Range("H136").Value = 4184.3335
$ws.Range("I136").Value = 2553
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 7659
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -2559
$ws.Range("N136").Value = -25200
$ws.Range("H137").Value = 3494
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 1816.6666
$ws.Range("J4").Value = 250
$ws.Range("L4").Value = 250
$ws.Range("N4").Value = -474
$ws.Range("H70").Value = 333333340
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 333333340
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H114").Value = 86662.336
$ws.Range("J114").Value = 86662.336
$ws.Range("L114").Value = 86662.336
$ws.Range("N114").Value = -95340.336
$ws.Range("H122").Value = 7814977
$ws.Range("I122").Value = 9616756
$ws.Range("K122").Value = 28850268
$ws.Range("M122").Value = -28847818
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 7750
$ws.Range("I20").Value = 7000
$ws.Range("K20").Value = 7000
$ws.Range("M20").Value = -6774
$ws.Range("H64").Value = 52500
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30450
$ws.Range("H67").Value = 52500
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31560
$ws.Range("H127").Value = 65355.75
$ws.Range("J127").Value = 65355.75
$ws.Range("L127").Value = 65355.75
$ws.Range("N127").Value = -75275.75
$ws.Range("H132").Value = 6192
$ws.Range("I132").Value = 5749.3335
$ws.Range("J132").Value = 6571.4287
$ws.Range("K132").Value = 17248.0005
$ws.Range("L132").Value = 19714.2861
$ws.Range("M132").Value = -14718.0005
$ws.Range("N132").Value = -24774.2861
$ws.Range("H136").Value = 2702.6667
$ws.Range("I136").Value = 2702.6667
$ws.Range("K136").Value = 8108.000100000001
$ws.Range("M136").Value = -5558.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H28").Value = 20019
$ws.Range("J28").Value = 20019
$ws.Range("L28").Value = 20019
$ws.Range("M28").Value = -20715
$ws.Range("H63").Value = 30000
$ws.Range("J63").Value = 30000
$ws.Range("L63").Value = 30000
$ws.Range("N63").Value = -31248
$ws.Range("H66").Value = 30000
$ws.Range("J66").Value = 30000
$ws.Range("L66").Value = 90000
$ws.Range("N66").Value = -96240
$ws.Range("H136").Value = 11159
$ws.Range("I136").Value = 10743.75
$ws.Range("J136").Value = 11989.5
$ws.Range("K136").Value = 32231.25
$ws.Range("L136").Value = 35968.5
$ws.Range("M136").Value = -29681.25
$ws.Range("N136").Value = -41068.5

Write-Host "Applied all Leve profit updates."